# menampilkan product di halaman product
# Update the "barang" (product) sheet with a refreshed product catalogue.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: A=id, B=id_category, C=nama_barang, D=slug, E=price, F=hpp,
#           G=deskripsi, H=size, I=qty
$data = @(
    @{ Row = 2;  A = 1;  B = 1; C = "Mistake";        D = "mistake";        E = 65000; F = 130000; G = "available"; H = "XL"; I = 50 },
    @{ Row = 3;  A = 2;  B = 1; C = "Basic Bk";        D = "basic-bk";       E = 65000; F = 130000; G = "available"; H = "L";  I = 50 },
    @{ Row = 4;  A = 3;  B = 1; C = "Licate";          D = "licate";         E = 65000; F = 130000; G = "available"; H = "M";  I = 50 },
    @{ Row = 5;  A = 4;  B = 2; C = "Anchor";          D = "anchor";         E = 90000; F = 185000; G = "available"; H = "M";  I = 50 },
    @{ Row = 6;  A = 5;  B = 2; C = "Crew";            D = "crew";           E = 90000; F = 185000; G = "available"; H = "L";  I = 50 },
    @{ Row = 7;  A = 6;  B = 3; C = "Daisy";           D = "daisy";          E = 75000; F = 145000; G = "available"; H = "XL"; I = 50 },
    @{ Row = 8;  A = 7;  B = 4; C = "Poppunk Suck";    D = "poppunk-suck";   E = 50000; F = 110000; G = "available"; H = "-";  I = 50 },
    @{ Row = 9;  A = 8;  B = 4; C = "Flag";            D = "Flag";           E = 50000; F = 110000; G = "available"; H = "-";  I = 50 },
    @{ Row = 10; A = 9;  B = 5; C = "Basic Gn";        D = "basic-gn";       E = 85000; F = 150000; G = "available"; H = "-";  I = 50 },
    @{ Row = 11; A = 10; B = 5; C = "Basic Logo";      D = "basic-logo";     E = 95000; F = 185000; G = "available"; H = "-";  I = 50 },
    @{ Row = 12; A = 11; B = 6; C = "Slop Anchor";     D = "slop-anchor";    E = 75000; F = 145000; G = "available"; H = "-";  I = 50 },
    @{ Row = 13; A = 12; B = 6; C = "Sunglass Basic";  D = "sunglass-basic"; E = 90000; F = 210000; G = "available"; H = "-";  I = 50 },
    @{ Row = 14; A = 13; B = 6; C = "Trapped";         D = "trapped";        E = 35000; F = 80000;  G = "available"; H = "-";  I = 50 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 5).NumberFormat = "#,##0"
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 6).NumberFormat = "#,##0"
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
    $ws.Cells.Item($r, 9).Value = $item.I
}

# Update the active selection to match the saved workbook state.
$ws.Range("K10").Select()
